$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "43.770.52"
$ws.Cells.Item(2, 5).Value = "  +0.42%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.292.40"
$ws.Cells.Item(3, 5).Value = "  +0.04%  "
$ws.Cells.Item(4, 5).Value = "  +0.17%  "
Set-TextValue $ws.Cells.Item(5, 4) "114.52"
$ws.Cells.Item(5, 5).Value = "  +19.03%  "
Set-TextValue $ws.Cells.Item(6, 4) "268.61"
$ws.Cells.Item(6, 5).Value = "  +0.40%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.616"
$ws.Cells.Item(7, 5).Value = "  -1.08%  "
$ws.Cells.Item(8, 5).Value = "  +0.21%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.623"
$ws.Cells.Item(9, 5).Value = "  +2.05%  "
Set-TextValue $ws.Cells.Item(10, 4) "48.17"
$ws.Cells.Item(10, 5).Value = "  +5.20%  "
$ws.Cells.Item(11, 5).Value = "  +0.16%  "
$ws.Cells.Item(12, 5).Value = "  +12.27%  "
Set-TextValue $ws.Cells.Item(13, 4) "0.107"
$ws.Cells.Item(13, 5).Value = "  +1.11%  "
$ws.Cells.Item(14, 5).Value = "  +2.99%  "
Set-TextValue $ws.Cells.Item(15, 4) "2.635.54"
$ws.Cells.Item(15, 5).Value = "  +0.16%  "
$ws.Cells.Item(16, 5).Value = "  -0.27%  "
Set-TextValue $ws.Cells.Item(17, 4) "2.292.23"
$ws.Cells.Item(17, 5).Value = "  +0.17%  "
Set-TextValue $ws.Cells.Item(18, 4) "43.609.90"
$ws.Cells.Item(18, 5).Value = "  +0.09%  "
$ws.Cells.Item(19, 5).Value = "  +2.22%  "
Set-TextValue $ws.Cells.Item(20, 4) "6.50"
$ws.Cells.Item(20, 5).Value = "  +4.78%  "
Set-TextValue $ws.Cells.Item(21, 4) "72.54"
$ws.Cells.Item(21, 5).Value = "  +0.38%  "
Set-TextValue $ws.Cells.Item(22, 4) "2.47"
$ws.Cells.Item(22, 5).Value = "  -2.60%  "
Set-TextValue $ws.Cells.Item(23, 4) "233.04"
$ws.Cells.Item(23, 5).Value = "  +0.00%  "
Set-TextValue $ws.Cells.Item(24, 4) "9.77"
$ws.Cells.Item(24, 5).Value = "  +6.68%  "
$ws.Cells.Item(25, 5).Value = "  +12.46%  "
$ws.Cells.Item(26, 5).Value = "  +0.00%  "
Set-TextValue $ws.Cells.Item(27, 4) "11.69"
$ws.Cells.Item(27, 5).Value = "  +5.18%  "
Set-TextValue $ws.Cells.Item(28, 4) "41.99"
$ws.Cells.Item(28, 5).Value = "  +3.54%  "
$ws.Cells.Item(29, 5).Value = "  -2.19%  "
$ws.Cells.Item(30, 5).Value = "  -0.16%  "
Set-TextValue $ws.Cells.Item(31, 4) "176.54"
$ws.Cells.Item(31, 5).Value = "  +0.51%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.0932"
$ws.Cells.Item(32, 5).Value = "  +4.43%  "
Set-TextValue $ws.Cells.Item(33, 4) "21.57"
$ws.Cells.Item(33, 5).Value = "  -1.15%  "
Set-TextValue $ws.Cells.Item(34, 4) "5.57"
$ws.Cells.Item(34, 5).Value = "  +3.87%  "
$ws.Cells.Item(35, 5).Value = "  +0.80%  "
Set-TextValue $ws.Cells.Item(36, 4) "4.75"
$ws.Cells.Item(36, 5).Value = "  +9.22%  "
$ws.Cells.Item(37, 5).Value = "  +0.56%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.0358"
$ws.Cells.Item(38, 5).Value = "  +1.12%  "
Set-TextValue $ws.Cells.Item(39, 4) "3.84"
$ws.Cells.Item(39, 5).Value = "  +13.16%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.45"
$ws.Cells.Item(40, 5).Value = "  +5.58%  "
$ws.Cells.Item(41, 5).Value = "  +13.07%  "
$ws.Cells.Item(42, 5).Value = "  +2.91%  "
Set-TextValue $ws.Cells.Item(43, 4) "72.81"
$ws.Cells.Item(43, 5).Value = "  +10.80%  "
$ws.Cells.Item(44, 2).Value = "THORChain"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Cells.Item(44, 4) "6.12"
$ws.Cells.Item(44, 5).Value = "  +17.23%  "

$ws.Cells.Item(45, 2).Value = "ARBITRUM"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Cells.Item(45, 4) "1.44"
$ws.Cells.Item(45, 5).Value = "  +6.41%  "

$ws.Cells.Item(46, 5).Value = "  +0.00%  "
Set-TextValue $ws.Cells.Item(47, 4) "8.71"
$ws.Cells.Item(47, 5).Value = "  -0.53%  "
Set-TextValue $ws.Cells.Item(48, 4) "102.98"
$ws.Cells.Item(48, 5).Value = "  +5.81%  "
$ws.Cells.Item(49, 5).Value = "  -1.26%  "
$ws.Cells.Item(50, 5).Value = "  +3.20%  "
$ws.Cells.Item(51, 5).Value = "  +4.09%  "
